# Add five new "corporate account" API sheets (insertAccounts, updateAccounts,
# deleteAccounts, getAccountContacts, addAccountContactAffiliations), modeled
# after the existing "getAccounts" sheet (URL / Param / SchemaPath layout).

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("getAccounts")

function Add-ApiSheet {
    param([string]$SheetName, [string]$Url, $Param, [string]$SchemaPath, [bool]$AddHyperlink)

    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $last)
    $ws.Name = $SheetName

    $ws.Range("A1").Value = "URL"
    $ws.Range("B1").Value = "Param"
    $ws.Range("C1").Value = "SchemaPath"
    $ws.Range("A2").Value = $Url
    $ws.Range("B2").Value = $Param
    $ws.Range("C2").Value = $SchemaPath

    if ($AddHyperlink) {
        $ws.Hyperlinks.Add($ws.Range("A2"), $Url) | Out-Null
    }

    # Match the look & feel of the other API sheets (bold header row, blue
    # underlined hyperlink cell) by copying formats from the template sheet.
    foreach ($cellRef in @("A1", "B1", "C1", "A2", "B2", "C2")) {
        $template.Range($cellRef).Copy() | Out-Null
        $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    }

    return $ws
}

Add-ApiSheet "insertAccounts" `
    "http://10.0.0.111:8080/web-corporate/api/submit/corporate/account/insertAccounts" `
    123 `
    "jsonSchema/corporate/account/insertAccounts.json" `
    $false | Out-Null

Add-ApiSheet "updateAccounts" `
    "http://localhost:8080/web-corporate/api/submit/corporate/account/updateAccounts" `
    '{"Header":{},"Params":{"Accounts":[{"AccountID":38,"oca":1}]}}' `
    "jsonSchema/corporate/account/updateAccounts.json" `
    $true | Out-Null

Add-ApiSheet "deleteAccounts" `
    "http://10.0.0.111:8080/web-corporate/api/submit/corporate/account/deleteAccounts" `
    '{"Header":{},"Params":{"AccountIDs": [43}}' `
    "jsonSchema/corporate/account/deleteAccounts.json" `
    $false | Out-Null

Add-ApiSheet "getAccountContacts" `
    "http://10.0.0.111:8080/web-corporate/api/submit/corporate/account/getAccountContacts" `
    '{"Header":{"PageSize":1000,"StartPosition":0,"Sort":[{"By":"ID","Dir":"DESC"}]},"Params":{"AcountIDs":[40]}}' `
    "jsonSchema/corporate/account/getAccountContacts.json" `
    $true | Out-Null

$lastWs = Add-ApiSheet "addAccountContactAffiliations" `
    "http://10.0.0.111:8080/web-corporate/api/submit/corporate/account/addAccountContactAffiliations" `
    '{"Header":{},"Params":{"AccountIDs":[39,40],"ContactIDs":[1,2],"AffiliationTypeIDs":[1,2]}}' `
    "jsonSchema/corporate/account/addAccountContactAffiliations.json" `
    $false

$lastWs.Select()

Write-Host "Sheets now:" $wb.Worksheets.Count
